$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Create the four new empty paragraphs up front (before any character
# formatting such as Bold is ever touched) so that no "current typing
# format" can leak from one paragraph into the next.
# ------------------------------------------------------------------
$pLast = $d.Paragraphs.Last
$null = $pLast.Range.InsertParagraphAfter()
$pE = $d.Paragraphs.Last
$pE.Style = "BodyText"

$null = $pE.Range.InsertParagraphAfter()
$pH = $d.Paragraphs.Last
$pH.Style = "Heading2"

$null = $pH.Range.InsertParagraphAfter()
$pT3 = $d.Paragraphs.Last
$pT3.Style = "FirstParagraph"

$null = $pT3.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Last
$pB.Style = "BodyText"

# ------------------------------------------------------------------
# 1) BodyText paragraph appended to the end of the "Text 2: E"
#    section, grading commentary that ends in a bold "E".
# ------------------------------------------------------------------
$r = $pE.Range
$r.Collapse(1)
$r.InsertAfter("Generellt kan eleven förbättra mycket i sin text, och därför sätter jag ett")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("E")
$r.Bold = 1
$r.Collapse(0)
$r.InsertAfter(".")

# ------------------------------------------------------------------
# 2) New "Text 3: A" Heading2 paragraph, bookmarked "text-3-a".
# ------------------------------------------------------------------
$r = $pH.Range
$r.Collapse(1)
$r.InsertAfter("Text 3: A")
$bmRange = $d.Range($pH.Range.Start, $pH.Range.End - 1)
$d.Bookmarks.Add("text-3-a", $bmRange)

# ------------------------------------------------------------------
# 3) New FirstParagraph paragraph discussing Text 3.
# ------------------------------------------------------------------
$r = $pT3.Range
$r.Collapse(1)
$r.InsertAfter("Text 3 verkar visa kvalite liknande text 1 på A-nivå. Övrigt tror jag att text 1 kör fram sina punkter lite bättre än denna text, men de är både på A-nivå. Texten är väl strukturerad, språket verkar lite mindre varierat (t.ex. börjar de flesta stycken med")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter([char]0x2018)
$r.Collapse(0)
$r.InsertAfter("jag")
$r.Collapse(0)
$r.InsertAfter([char]0x2019)
$r.Collapse(0)
$r.InsertAfter("), men ordval annars verkar rätt väl passande. Vissa delar av texten presenterar sig som mindre formella än några andra.")

# ------------------------------------------------------------------
# 4) New BodyText paragraph wrapping up with a bold "B".
# ------------------------------------------------------------------
$r = $pB.Range
$r.Collapse(1)
$r.InsertAfter("Jag skulle säga att texten befinner sig på en lägre A-nivå, mer likt ett")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("B")
$r.Bold = 1
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("där man får använda (d.v.s. inte här).")

Write-Output "edit applied"
